$wb = $excel.ActiveWorkbook

# Delete the row in "Inhalt" that describes the Ja_Nein worksheet mapping (row 11)
$inhalt = $wb.Worksheets.Item("Inhalt")
$inhalt.Rows.Item(11).Delete()

# Delete the "Ja_Nein" worksheet entirely
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Ja_Nein").Delete()
$excel.DisplayAlerts = $true

# Make sure "Inhalt" ends up the active/selected sheet
$inhalt.Activate()
$inhalt.Select()
